# Weekly data update: a new price record (week of 2021-12-15) is inserted
# as row 10, pushing the existing rows 10-113 down to rows 11-114 (so the
# data that used to be the last row, 113, becomes row 114).
#
# The new row 10 keeps the same market/category metadata (columns
# A,B,C,E,F,G,H,I,N,O,Q,R) as the row that ends up right below it (old row
# 10, now row 11) and only carries fresh values for the date/volume/price
# columns (D,J,K,L,M,P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10; rows 10..113 shift down to 11..114.
$ws.Rows.Item(10).Insert()

# Carry over the unchanged metadata columns from the row now just below.
$ws.Range("A10:C10").Value2 = $ws.Range("A11:C11").Value2
$ws.Range("E10:I10").Value2 = $ws.Range("E11:I11").Value2
$ws.Range("N10:O10").Value2 = $ws.Range("N11:O11").Value2
$ws.Range("Q10:R10").Value2 = $ws.Range("Q11:R11").Value2

# Fill in the new record's own data.
$ws.Range("D10").Value2 = 44545
$ws.Range("J10").Value2 = 80
$ws.Range("K10").Value2 = 8000
$ws.Range("L10").Value2 = 8500
$ws.Range("M10").Value2 = 8250
$ws.Range("P10").Value2 = 330
